$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap the full contents of row 10 and row 11 (data rows got
#     re-sorted upstream; only the row numbers stay put). Columns Y and AA
#     hold identical text-dates in both rows, so they're left untouched to
#     avoid Excel coercing the round-tripped string into a date serial.
$swapCols = @(
    "A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R",
    "S","T","U","V","W","X","Z","AB","AC","AD","AE","AF","AG","AH","AI","AJ",
    "AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY"
)

foreach ($c in $swapCols) {
    $addr10 = $c + "10"
    $addr11 = $c + "11"
    $v10 = $ws.Range($addr10).Value2
    $v11 = $ws.Range($addr11).Value2
    $ws.Range($addr10).Value = $v11
    $ws.Range($addr11).Value = $v10
}

# --- Step 2: every record's Taxonsorteringsordning (column B) shifts by +4
#     across all data rows (2-27).
for ($r = 2; $r -le 27; $r++) {
    $addr = "B" + $r
    $v = $ws.Range($addr).Value2
    $ws.Range($addr).Value = $v + 4
}
